# Weekly update for the "Hortaliza, Vega Monumental Concepción - Cilantro" sheet.
#
# A new week of data (report date 2021-09-14, Excel serial 44453) is inserted
# as a new pair of rows (Primera / Segunda) right after the pair currently
# sitting at rows 24:25 (2020-12-01), i.e. at rows 26:27. Every row from the
# old 26 onward shifts down by two rows (one pair), so the table grows from
# A1:R123 to A1:R125.
#
# The new pair is built by duplicating the row immediately above it (which
# carries all the constant "Cilantro" attributes and the standard
# Primera/Segunda price pattern) and then overwriting just the date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 26:27 - everything at/after row 26 shifts down by 2.
$ws.Rows("26:27").Insert()

# Populate the new Primera row (26) from the template directly above it (24),
# then the new Segunda row (27) from its template (25).
$ws.Range("A24:R24").Copy()
$ws.Range("A26:R26").PasteSpecial()

$ws.Range("A25:R25").Copy()
$ws.Range("A27:R27").PasteSpecial()

# Stamp the new pair with this week's report date (2021-09-14 -> serial 44453).
$ws.Range("D26").Value = 44453
$ws.Range("D27").Value = 44453
